$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.834.41'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '1.757.84'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''328.12'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '''0.4569'
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").Value = '''0.3496'
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("D9").Value = '''42.08'
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").Value = '''0.07346'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").Value = '''1.085'
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '''1.002'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '''20.60'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '''5.976'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = '''7.163'
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").Value = '1.758.14'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '''91.68'
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("D18").Value = '''0.00001052'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").Value = '''0.06420'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  -1.48%  '
$ws.Range("D22").Value = '''5.730'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").Value = '27.861.36'
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").Value = '''11.16'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = '''2.155'
$ws.Range("E25").Value = '  +3.45%  '
$ws.Range("D26").Value = '''162.11'
$ws.Range("E26").Value = '  -2.05%  '
$ws.Range("D27").Value = '''19.99'
$ws.Range("D28").Value = '1.960.02'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").Value = '''2.156'
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("D30").Value = '''123.12'
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("D31").Value = '''1.089'
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").Value = '''0.09305'
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").Value = '''3.641'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").Value = '''5.531'
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("D35").Value = '''11.75'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = '''0.06093'
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("D37").Value = '''0.02250'
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("D38").Value = '''0.2067'
$ws.Range("E38").Value = '  -1.04%  '
$ws.Range("D39").Value = '''4.893'
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("D40").Value = '''0.6185'
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("D41").Value = '''1.180'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '''1.367'
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("D43").Value = '''7.758'
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").Value = '''13.09'
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").Value = '''3.721'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").Value = '''0.5795'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '''122.04'
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").Value = '''1.922'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").Value = '''1.121'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").Value = '''0.06778'
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").Value = '''72.16'
$ws.Range("E51").Value = '  +0.39%  '
